$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.985.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +5.39%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.290.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.28%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'231.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.72%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.629"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.12%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'63.35"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.59%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.20%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.424"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.53%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0954"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +5.08%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'57.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.12%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'26.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +13.99%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.14%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.631.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.34%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.68%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'5.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +5.73%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.818"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.85%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.292.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.34%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'43.863.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +5.32%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +4.52%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'73.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.11%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.72%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'254.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.88%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'2.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +11.49%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.04%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.27%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.55%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'170.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.88%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.140"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.42%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'20.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.37%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.40%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +5.95%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.04%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.0704"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +8.63%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.24%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.46%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.38%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'6.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.22%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.25%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +3.72%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.01%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'8.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.42%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "'11.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +26.58%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "TerraClassic"
$ws.Range("C44").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D44").Value = "'0.000225"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.45%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.71%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'1.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.26%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'99.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.40%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0967"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.23%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'17.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.70%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.485.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.65%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'2.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.01%  "
$ws.Range("E51").Style = "Normal"
